# Insert a new "Extra Notes" slide as the 10th slide (pushing the existing
# 10th slide - "and when in doubt, ask the HelpDesk..." - down to position 11).
# Uses slide layout 3, the "Blank" layout that carries a Title placeholder
# and a Body (idx=1) placeholder - the same layout used by most of the
# other content slides in this deck.

$p = $ppt.ActivePresentation

$newSlide = $p.Slides.Add(10, 3)

$body = $newSlide.Shapes.Item(1)
$title = $newSlide.Shapes.Item(2)

$body.Name = "Text Placeholder 1"
$title.Name = "Title 2"

# ---- Title -----------------------------------------------------------
$title.TextFrame.TextRange.Text = "Extra Notes"

# ---- Body content ------------------------------------------------------
$paraTexts = @(
    "Cycle 2 proposal deadline: January 27th @ 8 PM EST = January 28th @ 12 PM AEDT.",
    "Proposals have total page limits but there are no limits on specific sections (i.e., Scientific Justification and Technical Justification sections are not limited in how long they can be within the total, but proposers should aim for a balance).",
    "A new version of the GTVT/MTVT (v0.3.0) was recently released. This release includes an updated JWST ephemeris, so updating this tool (if you have already installed it) would be a good idea to ensure more accurate visibility windows.",
    "JIST has not yet been updated based on commissioning results (waiting for release of ETC 2.0).",
    "The ""awesimsoss"" tool has been subsumed into MIRAGE and is no longer standalone.",
    "MIRISim can still be downloaded but is no longer under active development. Support for this tool is extremely reduced.",
    "The ""duty cycle"" definition is incorrect. It should be (saturation_time * number_integrations_per_exposure)/exposure_time.",
    "Simultaneous NIRCam coronagraphic imaging with the short- and long-wavelength channels is available in Cycle 2. APT and ETC will provide support for this updated mode.",
    "A new constraint, called the ""micrometeoroid avoidance zone (MAZ),"" will limit scheduling windows starting in Cycle 2. Basically, the number of observations in the direction of JWST's motion will be limited to reduce mirror degradation by micrometeoroid impacts. A new JDOX page will be published next week giving more details."
)

$bodyTF = $body.TextFrame
$bodyTF.TextRange.Text = [string]::Join([char]13, $paraTexts)
$bodyTF.AutoSize = 2

for ($i = 1; $i -le $paraTexts.Count; $i++) {
    $para = $bodyTF.TextRange.Paragraphs($i, 1)
    $para.ParagraphFormat.Alignment = 1
    $para.ParagraphFormat.Bullet.Visible = $true
    $para.ParagraphFormat.Bullet.Font.Name = "Arial"
    $para.ParagraphFormat.Bullet.Character = 8226
    $para.Font.Bold = $false
    $para.Font.Italic = $false
    $para.Font.Underline = $false
    $para.Font.Strikethrough = $false
    $para.Font.Shadow = $false
    $para.Font.Name = "Calibri"
    $para.Font.Color.RGB = 0
}

Write-Host "Inserted Extra Notes slide at index 10; total slides:" $p.Slides.Count
